# Update the two "taxtable" placeholder values with their real Dutch VAT
# labels (matches the SnelStart "Hoog"/"Laag" tax table names).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1").Value = "Hoog"
$ws.Range("R2").Value = "Laag"

# Move the selection/viewport to the last edited cell (R2), scrolled so
# that column M is the first visible column.
$ws.Range("R2").Select()
$excel.ActiveWindow.ScrollColumn = 13
